$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.871.95'
$ws.Range('E2').Value = '  -7.55%  '
$ws.Range('D3').Value = '2.525.59'
$ws.Range('E3').Value = '  -3.62%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = '296.64'
$ws.Range('E5').Value = '  -3.84%  '
$ws.Range('D6').Value = '93.45'
$ws.Range('E6').Value = '  -6.44%  '
$ws.Range('E7').Value = '  -5.22%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Value = '0.548'
$ws.Range('E9').Value = '  -5.81%  '
$ws.Range('D10').Value = '36.13'
$ws.Range('E10').Value = '  -8.32%  '
$ws.Range('E11').Value = '  -5.10%  '
$ws.Range('E12').Value = '  -6.36%  '
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('D14').Value = '2.910.49'
$ws.Range('E14').Value = '  -3.67%  '
$ws.Range('D15').Value = '2.509.24'
$ws.Range('E15').Value = '  -3.91%  '
$ws.Range('D16').Value = '0.865'
$ws.Range('E16').Value = '  -6.70%  '
$ws.Range('D17').Value = '14.13'
$ws.Range('E17').Value = '  -5.98%  '
$ws.Range('D18').Value = '42.862.56'
$ws.Range('E18').Value = '  -7.96%  '
$ws.Range('D19').Value = '6.61'
$ws.Range('E19').Value = '  -3.07%  '
$ws.Range('D20').Value = '0.0₃0963'
$ws.Range('E20').Value = '  -5.20%  '
$ws.Range('D21').Value = '12.24'
$ws.Range('E21').Value = '  -6.34%  '
$ws.Range('D22').Value = '72.69'
$ws.Range('E22').Value = '  +1.05%  '
$ws.Range('D23').Value = '259.52'
$ws.Range('E23').Value = '  -6.17%  '
$ws.Range('E24').Value = '  -4.70%  '
$ws.Range('D25').Value = '2.17'
$ws.Range('D26').Value = '29.02'
$ws.Range('E26').Value = '  -1.58%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').Value = '9.97'
$ws.Range('E28').Value = '  -6.55%  '
$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').Value = '36.78'
$ws.Range('E29').Value = '  -4.81%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '2.12'
$ws.Range('E30').Value = '  -7.12%  '
$ws.Range('D31').Value = '5.98'
$ws.Range('E31').Value = '  -7.45%  '
$ws.Range('D32').Value = '3.48'
$ws.Range('E32').Value = '  -4.39%  '
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').Value = '2.20'
$ws.Range('E33').Value = '  -2.88%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').Value = '151.47'
$ws.Range('E34').Value = '  -0.52%  '
$ws.Range('D35').Value = '2.74'
$ws.Range('E35').Value = '  -3.48%  '
$ws.Range('D36').Value = '0.0800'
$ws.Range('E36').Value = '  -5.05%  '
$ws.Range('D37').Value = '0.116'
$ws.Range('E37').Value = '  -6.11%  '
$ws.Range('E38').Value = '  -3.70%  '
$ws.Range('D39').Value = '23.79'
$ws.Range('E39').Value = '  -2.29%  '
$ws.Range('D40').Value = '16.44'
$ws.Range('E40').Value = '  +2.59%  '
$ws.Range('E41').Value = '  -4.60%  '
$ws.Range('E42').Value = '  -6.68%  '
$ws.Range('D43').Value = '3.83'
$ws.Range('E43').Value = '  -6.10%  '
$ws.Range('D44').Value = '2.024.72'
$ws.Range('E45').Value = '  -0.09%  '
$ws.Range('D46').Value = '85.69'
$ws.Range('E46').Value = '  -10.14%  '
$ws.Range('D47').Value = '1.60'
$ws.Range('D48').Value = '8.89'
$ws.Range('E48').Value = '  -6.40%  '
$ws.Range('D49').Value = '2.766.73'
$ws.Range('E49').Value = '  -3.85%  '
$ws.Range('D50').Value = '103.09'
$ws.Range('E50').Value = '  -6.04%  '
$ws.Range('E51').Value = '  -7.49%  '
